$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2: Id 13264 -> 34288
$ws.Range("A2").Value = 34288

# I2: Antal "2" -> blank
$ws.Range("I2").ClearContents()

# K2: Ålder-Stadium "larv" -> blank (cell removed)
$ws.Range("K2").ClearContents()

# P2: Lokalnamn "90-044, Venestad, Sk" -> "Venestad, Sk"
$ws.Range("P2").Value = "Venestad, Sk"

# S2: Noggrannhet 50 -> 25
$ws.Range("S2").Value = 25

# Y2: Startdatum "2008-09-01" -> "2008-03-01" (keep as text, not a date)
$ws.Range("Y2").NumberFormat = "@"
$ws.Range("Y2").Value = "2008-03-01"

# AA2: Slutdatum "2008-09-01" -> "2008-03-01" (keep as text, not a date)
$ws.Range("AA2").NumberFormat = "@"
$ws.Range("AA2").Value = "2008-03-01"

# AC2: Publik kommentar - new cell
$ws.Range("AC2").Value = "glest signalkräftbestånd, Mindre vattensalamander"

# AW2: Rapportör "Marika Stenberg" -> "Anders Hallengren"
$ws.Range("AW2").Value = "Anders Hallengren"

# AX2: Observatörer "Marika Stenberg, Per Nyström" -> "Per Nyström, Marika Stenberg"
$ws.Range("AX2").Value = "Per Nyström, Marika Stenberg"
